$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.731.36'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '2.251.73'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = "'307.71"
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = "'96.05"
$ws.Range("E6").Value = '  +2.25%  '
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").Value = "'35.56"
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("D11").Value = "'0.0807"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = "'7.27"
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").Value = '2.594.82'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.305.69'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = "'0.841"
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = "'13.64"
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '44.509.02'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'6.33"
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = "'12.06"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = "'65.71"
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").Value = "'238.90"
$ws.Range("E24").Value = '  +2.84%  '
$ws.Range("D25").Value = "'2.01"
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'2.28"
$ws.Range("E27").Value = '  +3.34%  '
$ws.Range("D28").Value = "'9.86"
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").Value = "'37.63"
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").Value = "'6.05"
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("D31").Value = "'20.01"
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").Value = "'152.68"
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("D33").Value = "'0.0800"
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").Value = "'2.64"
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("D35").Value = "'3.07"
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = "'1.86"
$ws.Range("E38").Value = '  +5.71%  '
$ws.Range("D39").Value = "'15.04"
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").Value = "'3.44"
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = "'3.80"
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("D42").Value = "'0.0302"
$ws.Range("E42").Value = '  +2.27%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = '1.831.63'
$ws.Range("E44").Value = '  +5.70%  '
$ws.Range("E45").Value = '  +17.56%  '
$ws.Range("D46").Value = "'0.193"
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("D47").Value = "'79.69"
$ws.Range("E47").Value = '  -5.76%  '
$ws.Range("D48").Value = "'71.05"
$ws.Range("E48").Value = '  +3.84%  '
$ws.Range("D49").Value = "'99.56"
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").Value = "'55.01"
$ws.Range("E51").Value = '  +2.48%  '
